# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# Updates the "Periodo Mora" labels (col E) and "Valor Mora" amounts (col F)
# for rows 16-53 of the account-statement table on Hoja1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 16; Periodo = "1603"; Valor = 27600 },
    @{ Row = 17; Periodo = "1604"; Valor = 27600 },
    @{ Row = 18; Periodo = "1605"; Valor = 27600 },
    @{ Row = 19; Periodo = "1606"; Valor = 27600 },
    @{ Row = 20; Periodo = "1607"; Valor = 25774 },
    @{ Row = 21; Periodo = "1608"; Valor = 25774 },
    @{ Row = 22; Periodo = "1609"; Valor = 25774 },
    @{ Row = 23; Periodo = "1610"; Valor = 25774 },
    @{ Row = 24; Periodo = "1611"; Valor = 25774 },
    @{ Row = 25; Periodo = "1612"; Valor = 25774 },
    @{ Row = 26; Periodo = "1701"; Valor = 25774 },
    @{ Row = 27; Periodo = "1702"; Valor = 25774 },
    @{ Row = 28; Periodo = "1703"; Valor = 25774 },
    @{ Row = 29; Periodo = "1704"; Valor = 25774 },
    @{ Row = 30; Periodo = "1705"; Valor = 25774 },
    @{ Row = 31; Periodo = "1706"; Valor = 25774 },
    @{ Row = 32; Periodo = "1707"; Valor = 25774 },
    @{ Row = 33; Periodo = "1708"; Valor = 25774 },
    @{ Row = 34; Periodo = "1709"; Valor = 25774 },
    @{ Row = 35; Periodo = "1710"; Valor = 25774 },
    @{ Row = 36; Periodo = "1810"; Valor = 31249 },
    @{ Row = 37; Periodo = "1811"; Valor = 31249 },
    @{ Row = 38; Periodo = "1812"; Valor = 31249 },
    @{ Row = 39; Periodo = "1901"; Valor = 31249 },
    @{ Row = 40; Periodo = "1902"; Valor = 31249 },
    @{ Row = 41; Periodo = "1903"; Valor = 31249 },
    @{ Row = 42; Periodo = "1904"; Valor = 31249 },
    @{ Row = 43; Periodo = "1905"; Valor = 31249 },
    @{ Row = 44; Periodo = "1906"; Valor = 31249 },
    @{ Row = 45; Periodo = "1907"; Valor = 31249 },
    @{ Row = 46; Periodo = "1908"; Valor = 31249 },
    @{ Row = 47; Periodo = "1909"; Valor = 31249 },
    @{ Row = 48; Periodo = "1910"; Valor = 31249 },
    @{ Row = 49; Periodo = "1911"; Valor = 31249 },
    @{ Row = 50; Periodo = "1912"; Valor = 31249 },
    @{ Row = 51; Periodo = "2001"; Valor = 31249 },
    @{ Row = 52; Periodo = "2002"; Valor = 31249 },
    @{ Row = 53; Periodo = "2003"; Valor = 31249 }
)

foreach ($item in $data) {
    $ws.Range("E$($item.Row)").Value = $item.Periodo
    $ws.Range("F$($item.Row)").Value = $item.Valor
}
